$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B:G, rows 2-8, per regenerated s_vals data
# (filters save games per commit message)

$data = @{
    2 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 5.582307763322248)
    3 = @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 1, 4.327115817150455)
    4 = @(0.6545652718822623, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 2.964545797025059)
    5 = @(0.2881169905109251, 9.983522426115931, 3.223369029078222, 13.86384647080068, 1, 27.35885491650576)
    6 = @(0.2881169905109251, 9.983522426115931, 0.7210945179870265, 13.86384647080068, 1, 24.85658040541457)
    7 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 6.15379541431027)
    8 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 6.15379541431027)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E
    $ws.Cells.Item($row, 6).Value = $vals[4]  # F
    $ws.Cells.Item($row, 7).Value = $vals[5]  # G
}
